# Refresh the cryptocurrency price/volume snapshot (and one Coin/Link swap)
# Matches the GitHub Actions scheduled data refresh for cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.029.43"
$ws.Range("E2").Value = "  +1.50%  "
# Row 3
$ws.Range("D3").Value = "1.767.36"
$ws.Range("E3").Value = "  -0.20%  "
# Row 4
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  +0.15%  "
# Row 5
$ws.Range("D5").Value = "'322.62"
$ws.Range("E5").Value = "  -0.57%  "
# Row 6
$ws.Range("D6").Value = "'0.9985"
$ws.Range("E6").Value = "  +0.01%  "
# Row 7
$ws.Range("D7").Value = "'0.4251"
$ws.Range("E7").Value = "  -4.50%  "
# Row 8
$ws.Range("D8").Value = "'0.3605"
$ws.Range("E8").Value = "  -3.61%  "
# Row 9
$ws.Range("D9").Value = "'44.34"
$ws.Range("E9").Value = "  -1.29%  "
# Row 10
$ws.Range("D10").Value = "'0.07481"
$ws.Range("E10").Value = "  -3.40%  "
# Row 11
$ws.Range("D11").Value = "'1.109"
$ws.Range("E11").Value = "  -1.01%  "
# Row 12
$ws.Range("D12").Value = "'0.9975"
$ws.Range("E12").Value = "  +0.22%  "
# Row 13
$ws.Range("D13").Value = "'21.50"
$ws.Range("E13").Value = "  -1.02%  "
# Row 14
$ws.Range("D14").Value = "'6.110"
$ws.Range("E14").Value = "  -1.20%  "
# Row 15
$ws.Range("D15").Value = "'7.367"
$ws.Range("E15").Value = "  -1.01%  "
# Row 16
$ws.Range("D16").Value = "1.791.59"
$ws.Range("E16").Value = "  +1.59%  "
# Row 17
$ws.Range("D17").Value = "'91.81"
$ws.Range("E17").Value = "  +1.25%  "
# Row 18
$ws.Range("D18").Value = "'0.00001063"
$ws.Range("E18").Value = "  -1.14%  "
# Row 19
$ws.Range("D19").Value = "'0.06400"
$ws.Range("E19").Value = "  +2.10%  "
# Row 20
$ws.Range("D20").Value = "'0.9976"
$ws.Range("E20").Value = "  -0.03%  "
# Row 21
$ws.Range("D21").Value = "'17.16"
$ws.Range("E21").Value = "  -1.62%  "
# Row 22
$ws.Range("D22").Value = "'5.992"
$ws.Range("E22").Value = "  -3.37%  "
# Row 23
$ws.Range("D23").Value = "28.020.43"
$ws.Range("E23").Value = "  +1.44%  "
# Row 24
$ws.Range("E24").Value = "  -2.24%  "
# Row 25
$ws.Range("D25").Value = "'2.142"
$ws.Range("E25").Value = "  -7.39%  "
# Row 26
$ws.Range("D26").Value = "'158.66"
# Row 27
$ws.Range("D27").Value = "'20.22"
$ws.Range("E27").Value = "  -2.19%  "
# Row 28
$ws.Range("D28").Value = "1.990.48"
$ws.Range("E28").Value = "  +1.59%  "
# Row 29
$ws.Range("D29").Value = "'2.152"
$ws.Range("E29").Value = "  -7.10%  "
# Row 30
$ws.Range("D30").Value = "'126.14"
$ws.Range("E30").Value = "  -1.20%  "
# Row 31
$ws.Range("D31").Value = "'1.179"
$ws.Range("E31").Value = "  -0.28%  "
# Row 32
$ws.Range("D32").Value = "'5.685"
$ws.Range("E32").Value = "  -0.99%  "
# Row 33
$ws.Range("D33").Value = "'0.09050"
$ws.Range("E33").Value = "  -1.63%  "
# Row 34
$ws.Range("D34").Value = "'3.516"
$ws.Range("E34").Value = "  -3.10%  "
# Row 35
$ws.Range("D35").Value = "'12.62"
$ws.Range("E35").Value = "  -0.29%  "
# Row 36
$ws.Range("D36").Value = "'0.02337"
$ws.Range("E36").Value = "  +0.61%  "
# Row 37
$ws.Range("D37").Value = "'5.061"
$ws.Range("E37").Value = "  -0.11%  "
# Row 38
$ws.Range("D38").Value = "'0.2112"
$ws.Range("E38").Value = "  -2.40%  "
# Row 39
$ws.Range("D39").Value = "'0.06077"
$ws.Range("E39").Value = "  -1.01%  "
# Row 40
$ws.Range("E40").Value = "  -1.06%  "
# Row 41
$ws.Range("D41").Value = "'1.188"
$ws.Range("E41").Value = "  +0.51%  "
# Row 42
$ws.Range("D42").Value = "'0.9980"
$ws.Range("E42").Value = "  +0.07%  "
# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'7.828"
$ws.Range("E43").Value = "  -1.91%  "
# Row 44
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.392"
$ws.Range("E44").Value = "  +0.34%  "
# Row 45
$ws.Range("D45").Value = "'13.59"
$ws.Range("E45").Value = "  -0.35%  "
# Row 46
$ws.Range("D46").Value = "'0.5969"
$ws.Range("E46").Value = "  -0.44%  "
# Row 47
$ws.Range("E47").Value = "  -0.60%  "
# Row 48
$ws.Range("D48").Value = "'2.019"
$ws.Range("E48").Value = "  +1.54%  "
# Row 49
$ws.Range("D49").Value = "'123.54"
$ws.Range("E49").Value = "  -2.18%  "
# Row 50
$ws.Range("D50").Value = "'1.176"
$ws.Range("E50").Value = "  +3.41%  "
# Row 51
$ws.Range("D51").Value = "'0.06882"
$ws.Range("E51").Value = "  -0.01%  "
